# Upload Y4_B2526_Excuses.xlsx via attendance app
# Updates Student ID (col A) and Log Date (col C) for rows 2-7.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; StudentId = "200933"; LogDate = "25/10/2025" },
    @{ Row = 3; StudentId = "200792"; LogDate = "25/10/2025" },
    @{ Row = 4; StudentId = "211177"; LogDate = "25/10/2025" },
    @{ Row = 5; StudentId = "191088"; LogDate = "25/10/2025" },
    @{ Row = 6; StudentId = "211169"; LogDate = "25/10/2025" },
    @{ Row = 7; StudentId = "211741"; LogDate = "25/10/2025" }
)

foreach ($u in $updates) {
    $cellA = $ws.Cells.Item($u.Row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $u.StudentId

    $cellC = $ws.Cells.Item($u.Row, 3)
    $cellC.NumberFormat = "@"
    $cellC.Value = $u.LogDate
}
